$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36 (100k resistors): remove R56 from the reference-designator list ---
# Before: "R11,R14,R17,R20,R35,R36,R37,R38,R48,R49,R55,R56"
# After:  "R11,R14,R17,R20,R35,R36,R37,R38,R48,R49,R55"
$c36 = $ws.Range("C36")
$c36.Characters(44, 4).Text = ""

# --- Row 31 (1k resistors): remove R58 from the reference-designator list ---
# Before: "R10,R13,R16,R19,R21,R23,R24,R29,R30,R39,R40,R50,R51,R57,R58,R59,R60,R62,R64"
# After:  "R10,R13,R16,R19,R21,R23,R24,R29,R30,R39,R40,R50,R51,R57,R59,R60,R62,R64"
# This cell carries rich-text colouring (R39,R40,R59,R60 green; R64 red) that gets
# flattened to plain text by the Characters() text edit, so re-apply the colours
# to the surviving runs afterwards.
$c31 = $ws.Range("C31")
$c31.Characters(57, 4).Text = ""

$green = 5287936   # RGB(0,176,80)  == FF00B050
$red = 255         # RGB(255,0,0)   == FFFF0000
$black = 0         # RGB(0,0,0)     == FF000000

$c31.Characters(37, 3).Font.Color = $green   # R39
$c31.Characters(40, 1).Font.Color = $black   # ,
$c31.Characters(41, 3).Font.Color = $green   # R40
$c31.Characters(44, 13).Font.Color = $black  # ,R50,R51,R57,
$c31.Characters(57, 3).Font.Color = $green   # R59
$c31.Characters(60, 1).Font.Color = $black   # ,
$c31.Characters(61, 3).Font.Color = $green   # R60
$c31.Characters(64, 5).Font.Color = $black   # ,R62,
$c31.Characters(69, 3).Font.Color = $red     # R64
